$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken column with per-row timestamps
$ws.Range("F2").Value = "2021-10-05 10:50:05.642279"
$ws.Range("F3").Value = "2021-10-05 10:50:05.642291"
$ws.Range("F4").Value = "2021-10-05 10:50:05.642295"
$ws.Range("F5").Value = "2021-10-05 10:50:05.642297"
$ws.Range("F6").Value = "2021-10-05 10:50:05.642300"
$ws.Range("F7").Value = "2021-10-05 10:50:05.642303"
$ws.Range("F8").Value = "2021-10-05 10:50:05.642305"
$ws.Range("F9").Value = "2021-10-05 10:50:05.642308"
$ws.Range("F10").Value = "2021-10-05 10:50:05.642311"
$ws.Range("F11").Value = "2021-10-05 10:50:05.642313"
$ws.Range("F12").Value = "2021-10-05 10:50:05.642316"
$ws.Range("F13").Value = "2021-10-05 10:50:05.642318"
$ws.Range("F14").Value = "2021-10-05 10:50:05.642320"
$ws.Range("F15").Value = "2021-10-05 10:50:05.642323"
$ws.Range("F16").Value = "2021-10-05 10:50:05.642325"
$ws.Range("F17").Value = "2021-10-05 10:50:05.642328"
$ws.Range("F18").Value = "2021-10-05 10:50:05.642331"
$ws.Range("F19").Value = "2021-10-05 10:50:05.642333"
$ws.Range("F20").Value = "2021-10-05 10:50:05.642336"
$ws.Range("F21").Value = "2021-10-05 10:50:05.642338"
$ws.Range("F22").Value = "2021-10-05 10:50:05.642341"
$ws.Range("F23").Value = "2021-10-05 10:50:05.642344"
$ws.Range("F24").Value = "2021-10-05 10:50:05.642346"
$ws.Range("F25").Value = "2021-10-05 10:50:05.642349"
$ws.Range("F26").Value = "2021-10-05 10:50:05.642352"
$ws.Range("F27").Value = "2021-10-05 10:50:05.642354"
$ws.Range("F28").Value = "2021-10-05 10:50:05.642357"
$ws.Range("F29").Value = "2021-10-05 10:50:05.642359"
$ws.Range("F30").Value = "2021-10-05 10:50:05.642362"
$ws.Range("F31").Value = "2021-10-05 10:50:05.642364"
$ws.Range("F32").Value = "2021-10-05 10:50:05.642367"
$ws.Range("F33").Value = "2021-10-05 10:50:05.642369"
$ws.Range("F34").Value = "2021-10-05 10:50:05.642372"
$ws.Range("F35").Value = "2021-10-05 10:50:05.642375"
$ws.Range("F36").Value = "2021-10-05 10:50:05.642377"
$ws.Range("F37").Value = "2021-10-05 10:50:05.642380"
$ws.Range("F38").Value = "2021-10-05 10:50:05.642382"
$ws.Range("F39").Value = "2021-10-05 10:50:05.642385"
$ws.Range("F40").Value = "2021-10-05 10:50:05.642387"

Write-Output "done"
